$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 189, shifting row 189 and everything below it
# down by one (so the former row 189 becomes row 190, ..., former row 224
# becomes row 225).
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with its data.
$ws.Range("A189").Value = 8
$ws.Range("B189").Value = "Terminal La Palmera de La Serena"
$ws.Range("C189").Value = "Coquimbo"
$ws.Range("D189").Value = 44694
$ws.Range("E189").Value = 4
$ws.Range("F189").Value = 100112031
$ws.Range("G189").Value = "Poroto verde"
$ws.Range("H189").Value = "Magnum"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 480
$ws.Range("K189").Value = 26000
$ws.Range("L189").Value = 27000
$ws.Range("M189").Value = 26500
$ws.Range("N189").Value = "$/malla 25 kilos"
$ws.Range("O189").Value = "Provincia de Limarí"
$ws.Range("P189").Value = 1060
$ws.Range("Q189").Value = 25
$ws.Range("R189").Value = "Hortaliza"
